$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 833660.25
$ws.Range("J6").Value = 679.75
$ws.Range("L6").Value = 2039.25
$ws.Range("N6").Value = -2263.25

$ws.Range("H8").Value = 241.75
$ws.Range("I8").Value = 140.57143
$ws.Range("K8").Value = 421.71429
$ws.Range("M8").Value = -282.71429

$ws.Range("H86").Value = 11114667
$ws.Range("I86").Value = 1550.5834
$ws.Range("J86").Value = 33340900
$ws.Range("K86").Value = 1550.5834
$ws.Range("L86").Value = 33340900
$ws.Range("M86").Value = -427.5834
$ws.Range("N86").Value = -33343146

$ws.Range("H89").Value = 11114667
$ws.Range("I89").Value = 1550.5834
$ws.Range("J89").Value = 33340900
$ws.Range("K89").Value = 7752.916999999999
$ws.Range("L89").Value = 166704500
$ws.Range("M89").Value = -2136.916999999999
$ws.Range("N89").Value = -166715732

$ws.Range("H138").Value = 1446.5139
$ws.Range("I138").Value = 970.9400000000001
$ws.Range("J138").Value = 2527.3635
$ws.Range("K138").Value = 2912.82
$ws.Range("L138").Value = 7582.0905
$ws.Range("M138").Value = 2227.18
$ws.Range("N138").Value = -17862.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 3400.6
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 4000
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 4000
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -4346

$ws.Range("H17").Value = 5450
$ws.Range("J17").Value = 5450
$ws.Range("L17").Value = 5450
$ws.Range("N17").Value = -5796

$ws.Range("H32").Value = 20886.518
$ws.Range("I32").Value = 4726.7124
$ws.Range("J32").Value = 164529.22
$ws.Range("K32").Value = 4726.7124
$ws.Range("L32").Value = 164529.22
$ws.Range("M32").Value = -4439.7124
$ws.Range("N32").Value = -165103.22

$ws.Range("H61").Value = 1895.5625
$ws.Range("I61").Value = 998.3913
$ws.Range("J61").Value = 2720.96
$ws.Range("K61").Value = 998.3913
$ws.Range("L61").Value = 2720.96
$ws.Range("M61").Value = -786.3913
$ws.Range("N61").Value = -3144.96

$ws.Range("H110").Value = 52742696
$ws.Range("I110").Value = 77085040
$ws.Range("J110").Value = 950
$ws.Range("K110").Value = 77085040
$ws.Range("L110").Value = 950
$ws.Range("M110").Value = -77082995
$ws.Range("N110").Value = -5040

$ws.Range("H136").Value = 1895.5625
$ws.Range("I136").Value = 998.3913
$ws.Range("J136").Value = 2720.96
$ws.Range("K136").Value = 2995.1739
$ws.Range("L136").Value = 8162.88
$ws.Range("M136").Value = -445.1738999999998
$ws.Range("N136").Value = -13262.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 400
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 400
$ws.Range("N7").Value = -626
$ws.Range("M7").ClearContents()

$ws.Range("H8").Value = 26367.334
$ws.Range("I8").Value = 202
$ws.Range("J8").Value = 39450
$ws.Range("K8").Value = 202
$ws.Range("L8").Value = 39450
$ws.Range("M8").Value = -62
$ws.Range("N8").Value = -39730

$ws.Range("H10").Value = 10000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10280

$ws.Range("H11").Value = 856.6667
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 828
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 828
$ws.Range("M11").Value = -860
$ws.Range("N11").Value = -1108

$ws.Range("H12").Value = 996.25
$ws.Range("I12").Value = 2005
$ws.Range("J12").Value = 660
$ws.Range("K12").Value = 2005
$ws.Range("L12").Value = 660
$ws.Range("M12").Value = -1837
$ws.Range("N12").Value = -996

$ws.Range("H14").Value = 34980
$ws.Range("J14").Value = 34980
$ws.Range("L14").Value = 34980
$ws.Range("N14").Value = -35324

$ws.Range("H16").Value = 7102.3335
$ws.Range("I16").Value = 753.5
$ws.Range("J16").Value = 19800
$ws.Range("K16").Value = 753.5
$ws.Range("L16").Value = 19800
$ws.Range("M16").Value = -583.5
$ws.Range("N16").Value = -20140

$ws.Range("H17").Value = 5333.3335
$ws.Range("J17").Value = 5333.3335
$ws.Range("L17").Value = 5333.3335
$ws.Range("N17").Value = -5677.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.04000000000001
$ws.Range("J7").Value = 173.57143
$ws.Range("L7").Value = 173.57143
$ws.Range("N7").Value = -399.57143

$ws.Range("H19").Value = 24101.666
$ws.Range("I19").Value = 36.666668
$ws.Range("J19").Value = 48166.668
$ws.Range("K19").Value = 36.666668
$ws.Range("L19").Value = 48166.668
$ws.Range("M19").Value = 133.333332
$ws.Range("N19").Value = -48506.668

$ws.Range("H24").Value = 24101.666
$ws.Range("I24").Value = 36.666668
$ws.Range("J24").Value = 48166.668
$ws.Range("K24").Value = 36.666668
$ws.Range("L24").Value = 48166.668
$ws.Range("M24").Value = 133.333332
$ws.Range("N24").Value = -48506.668

$ws.Range("H58").Value = 568
$ws.Range("I58").Value = 522
$ws.Range("J58").Value = 614
$ws.Range("K58").Value = 522
$ws.Range("L58").Value = 614
$ws.Range("M58").Value = -319
$ws.Range("N58").Value = -1020

$ws.Range("H136").Value = 568
$ws.Range("I136").Value = 522
$ws.Range("J136").Value = 614
$ws.Range("K136").Value = 1566
$ws.Range("L136").Value = 1842
$ws.Range("M136").Value = 984
$ws.Range("N136").Value = -6942

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 6641.1113
$ws.Range("I121").Value = 5655.8887
$ws.Range("J121").Value = 7626.3335
$ws.Range("K121").Value = 16967.6661
$ws.Range("L121").Value = 22879.0005
$ws.Range("M121").Value = -15657.6661
$ws.Range("N121").Value = -25499.0005

$ws.Range("H132").Value = 1838.5
$ws.Range("I132").Value = 1081
$ws.Range("J132").Value = 2596
$ws.Range("K132").Value = 9729
$ws.Range("L132").Value = 23364
$ws.Range("M132").Value = -7199
$ws.Range("N132").Value = -28424

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 716
$ws.Range("I9").Value = 716
$ws.Range("K9").Value = 716
$ws.Range("M9").Value = -546

$ws.Range("H12").Value = 8000800.5
$ws.Range("I12").Value = 8000800.5
$ws.Range("K12").Value = 8000800.5
$ws.Range("M12").Value = -8000660.5

$ws.Range("H13").Value = 313.75
$ws.Range("I13").Value = 313.75
$ws.Range("K13").Value = 313.75
$ws.Range("M13").Value = -174.75

$ws.Range("H70").Value = 66826
$ws.Range("I70").Value = 121261.766
$ws.Range("J70").Value = 5132.1333
$ws.Range("K70").Value = 121261.766
$ws.Range("L70").Value = 5132.1333
$ws.Range("M70").Value = -120991.766
$ws.Range("N70").Value = -5672.1333

$ws.Range("H73").Value = 66826
$ws.Range("I73").Value = 121261.766
$ws.Range("J73").Value = 5132.1333
$ws.Range("K73").Value = 121261.766
$ws.Range("L73").Value = 5132.1333
$ws.Range("M73").Value = -120325.766
$ws.Range("N73").Value = -7004.1333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 21442
$ws.Range("I19").Value = 3136.6667
$ws.Range("J19").Value = 48900
$ws.Range("K19").Value = 3136.6667
$ws.Range("L19").Value = 48900
$ws.Range("M19").Value = -2966.6667
$ws.Range("N19").Value = -49240

$ws.Range("H108").Value = 26406
$ws.Range("J108").Value = 26406
$ws.Range("L108").Value = 26406
$ws.Range("N108").Value = -34086

$ws.Range("H132").Value = 3031.0322
$ws.Range("I132").Value = 3128.0344
$ws.Range("J132").Value = 1624.5
$ws.Range("K132").Value = 9384.1032
$ws.Range("L132").Value = 4873.5
$ws.Range("M132").Value = -6854.1032
$ws.Range("N132").Value = -9933.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 868.5
$ws.Range("I17").Value = 868.5
$ws.Range("K17").Value = 868.5
$ws.Range("M17").Value = -696.5

$ws.Range("H33").Value = 39915.25
$ws.Range("I33").Value = 9819
$ws.Range("J33").Value = 49947.332
$ws.Range("K33").Value = 9819
$ws.Range("L33").Value = 49947.332
$ws.Range("M33").Value = -9569
$ws.Range("N33").Value = -50447.332

$ws.Range("H36").Value = 39915.25
$ws.Range("I36").Value = 9819
$ws.Range("J36").Value = 49947.332
$ws.Range("K36").Value = 9819
$ws.Range("L36").Value = 49947.332
$ws.Range("M36").Value = -9569
$ws.Range("N36").Value = -50447.332

$ws.Range("H54").Value = 6795.1816
$ws.Range("J54").Value = 6767.7
$ws.Range("L54").Value = 6767.7
$ws.Range("N54").Value = -7807.7

$ws.Range("H81").Value = 223829
$ws.Range("I81").Value = 250875.5
$ws.Range("J81").Value = 202191.8
$ws.Range("K81").Value = 501751
$ws.Range("L81").Value = 404383.6
$ws.Range("M81").Value = -500690
$ws.Range("N81").Value = -406505.6

$ws.Range("H84").Value = 223829
$ws.Range("I84").Value = 250875.5
$ws.Range("J84").Value = 202191.8
$ws.Range("K84").Value = 2508755
$ws.Range("L84").Value = 2021918
$ws.Range("M84").Value = -2503451
$ws.Range("N84").Value = -2032526

$ws.Range("H132").Value = 2420.8462
$ws.Range("I132").Value = 2655.1614
$ws.Range("J132").Value = 1512.875
$ws.Range("K132").Value = 7965.4842
$ws.Range("L132").Value = 4538.625
$ws.Range("M132").Value = -5435.4842
$ws.Range("N132").Value = -9598.625
